$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("AppControl")
$ws1.Activate()
$ws1.Range("B25").Select()
$aw = $excel.ActiveWindow
try { $aw.ScrollRow = 25; "aw ScrollRow ok" } catch { "aw ScrollRow err: $_" }
try { $aw.ScrollColumn = 1; "aw ScrollColumn ok" } catch { "aw ScrollColumn err: $_" }
"current scrollrow: " + $aw.ScrollRow
